$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "Pre Experimental Phase" (column C) measurements for session 3
$ws.Range("C2").Value = 55
$ws.Range("C3").Value = 51
$ws.Range("C4").Value = 53
$ws.Range("C5").Value = 71
$ws.Range("C6").Value = 75
$ws.Range("C7").Value = 75
$ws.Range("C9").Value = 78
$ws.Range("C10").Value = 63
$ws.Range("C11").Value = 82
$ws.Range("C13").Value = 75
$ws.Range("C14").Value = 71

# Update the active selection to C15 (next empty row in column C)
$ws.Range("C15").Select()
